# Update gh-pages output (苏州-漫展信息.xlsx) to match the data generated
# at commit 456a3b4: a couple of "want to go" counters were refreshed, and
# a new performance ("苏州·春日计划2024——特别二次元不插电音乐会", on
# 2024-06-01) was published ahead of the existing 2024-06-02 entry, so it
# is inserted as a new row on the sheets that list it.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibition) - counters refreshed, no new rows here.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 14864
$wsExpo.Range("F3").Value = 18465
$wsExpo.Range("F5").Value = 111
$wsExpo.Range("F17").Value = 1407
$wsExpo.Range("F22").Value = 7657
$wsExpo.Range("F28").Value = 5949
$wsExpo.Range("F34").Value = 5290

# ---------------------------------------------------------------------------
# Sheet "演出" (Performance) - insert the new 2024-06-01 show as row 3,
# pushing the existing 2024-06-02 show down to row 4.
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Rows.Item(3).Insert()

$wsShow.Range("A3").Value = 2
$wsShow.Range("B3").NumberFormat = "@"
$wsShow.Range("B3").Value = "2024-06-01"
$wsShow.Range("C3").Value = "苏州·春日计划2024——特别二次元不插电音乐会"
$wsShow.Range("D3").Value = "星湖街555号高教区(体育馆南侧) 苏州独墅湖影剧院"
$wsShow.Range("E3").Value = "2024.06.01 19:30-06.01 21:00"
$wsShow.Range("F3").Value = 0
$wsShow.Range("G3").Value = 88
$wsShow.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=84720"
$wsShow.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202404/gwLWvSew1713796405109.png"

# Row that got pushed down keeps its content; only its running index moves.
$wsShow.Range("A4").Value = 3

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - mirrors both the counter refresh and the
# new performance row above (inserted at row 29 here).
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 14864
$wsAll.Range("F3").Value = 18465
$wsAll.Range("F5").Value = 111
$wsAll.Range("F17").Value = 1407
$wsAll.Range("F23").Value = 7657

$wsAll.Rows.Item(29).Insert()

$wsAll.Range("A29").Value = 28
$wsAll.Range("B29").NumberFormat = "@"
$wsAll.Range("B29").Value = "2024-06-01"
$wsAll.Range("C29").Value = "苏州·春日计划2024——特别二次元不插电音乐会"
$wsAll.Range("D29").Value = "星湖街555号高教区(体育馆南侧) 苏州独墅湖影剧院"
$wsAll.Range("E29").Value = "2024.06.01 19:30-06.01 21:00"
$wsAll.Range("F29").Value = 0
$wsAll.Range("G29").Value = 88
$wsAll.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=84720"
$wsAll.Range("I29").Value = "//i1.hdslb.com/bfs/openplatform/202404/gwLWvSew1713796405109.png"

# Re-sequence the index column for the rows that shifted down one place
# (was 28..37, now 29..38 for rows 30..39).
$wsAll.Range("A30").Value = 29
$wsAll.Range("A31").Value = 30
$wsAll.Range("A32").Value = 31
$wsAll.Range("A33").Value = 32
$wsAll.Range("A34").Value = 33
$wsAll.Range("A35").Value = 34
$wsAll.Range("A36").Value = 35
$wsAll.Range("A37").Value = 36
$wsAll.Range("A38").Value = 37
$wsAll.Range("A39").Value = 38

# Two more counters refreshed for events that now sit one row lower:
# row 31 = "Come in joy" (was row 30), row 37 = "萤火国潮文化节" (was row 36).
$wsAll.Range("F31").Value = 5949
$wsAll.Range("F37").Value = 5290
